$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy formats from the last existing row (506) down across the new rows (507:547)
$fmtSrc = $ws.Range("A506:I506")
$fmtDst = $ws.Range("A507:I547")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Write the values for columns A-F for each new row
$ws.Cells.Item(507,1).Value2 = "Pakistan"
$ws.Cells.Item(507,2).Value2 = 44391
$ws.Cells.Item(507,3).Value2 = 981392
$ws.Cells.Item(507,4).Value2 = 22689
$ws.Cells.Item(507,5).Value2 = 916373
$ws.Cells.Item(507,6).Value2 = 42330
$ws.Cells.Item(508,1).Value2 = "Pakistan"
$ws.Cells.Item(508,2).Value2 = 44392
$ws.Cells.Item(508,3).Value2 = 983719
$ws.Cells.Item(508,4).Value2 = 22720
$ws.Cells.Item(508,5).Value2 = 917329
$ws.Cells.Item(508,6).Value2 = 43670
$ws.Cells.Item(509,1).Value2 = "Pakistan"
$ws.Cells.Item(509,2).Value2 = 44393
$ws.Cells.Item(509,3).Value2 = 986668
$ws.Cells.Item(509,4).Value2 = 22760
$ws.Cells.Item(509,5).Value2 = 918329
$ws.Cells.Item(509,6).Value2 = 43670
$ws.Cells.Item(510,1).Value2 = "Pakistan"
$ws.Cells.Item(510,2).Value2 = 44394
$ws.Cells.Item(510,3).Value2 = 989275
$ws.Cells.Item(510,4).Value2 = 22781
$ws.Cells.Item(510,5).Value2 = 919163
$ws.Cells.Item(510,6).Value2 = 47331
$ws.Cells.Item(511,1).Value2 = "Pakistan"
$ws.Cells.Item(511,2).Value2 = 44395
$ws.Cells.Item(511,3).Value2 = 991727
$ws.Cells.Item(511,4).Value2 = 22811
$ws.Cells.Item(511,5).Value2 = 920066
$ws.Cells.Item(511,6).Value2 = 48850
$ws.Cells.Item(512,1).Value2 = "Pakistan"
$ws.Cells.Item(512,2).Value2 = 44396
$ws.Cells.Item(512,3).Value2 = 993872
$ws.Cells.Item(512,4).Value2 = 22848
$ws.Cells.Item(512,5).Value2 = 921095
$ws.Cells.Item(512,6).Value2 = 48850
$ws.Cells.Item(513,1).Value2 = "Pakistan"
$ws.Cells.Item(513,2).Value2 = 44397
$ws.Cells.Item(513,3).Value2 = 996451
$ws.Cells.Item(513,4).Value2 = 22888
$ws.Cells.Item(513,5).Value2 = 922034
$ws.Cells.Item(513,6).Value2 = 51529
$ws.Cells.Item(514,1).Value2 = "Pakistan"
$ws.Cells.Item(514,2).Value2 = 44398
$ws.Cells.Item(514,3).Value2 = 998609
$ws.Cells.Item(514,4).Value2 = 22928
$ws.Cells.Item(514,5).Value2 = 922929
$ws.Cells.Item(514,6).Value2 = 52752
$ws.Cells.Item(515,1).Value2 = "Pakistan"
$ws.Cells.Item(515,2).Value2 = 44399
$ws.Cells.Item(515,3).Value2 = 1000034
$ws.Cells.Item(515,4).Value2 = 22939
$ws.Cells.Item(515,5).Value2 = 923472
$ws.Cells.Item(515,6).Value2 = 53623
$ws.Cells.Item(516,1).Value2 = "Pakistan"
$ws.Cells.Item(516,2).Value2 = 44400
$ws.Cells.Item(516,3).Value2 = 1001875
$ws.Cells.Item(516,4).Value2 = 22971
$ws.Cells.Item(516,5).Value2 = 924782
$ws.Cells.Item(516,6).Value2 = 54122
$ws.Cells.Item(517,1).Value2 = "Pakistan"
$ws.Cells.Item(517,2).Value2 = 44401
$ws.Cells.Item(517,3).Value2 = 1004694
$ws.Cells.Item(517,4).Value2 = 23016
$ws.Cells.Item(517,5).Value2 = 925958
$ws.Cells.Item(517,6).Value2 = 55720
$ws.Cells.Item(518,1).Value2 = "Pakistan"
$ws.Cells.Item(518,2).Value2 = 44402
$ws.Cells.Item(518,3).Value2 = 1008446
$ws.Cells.Item(518,4).Value2 = 23048
$ws.Cells.Item(518,5).Value2 = 927599
$ws.Cells.Item(518,6).Value2 = 57799
$ws.Cells.Item(519,1).Value2 = "Pakistan"
$ws.Cells.Item(519,2).Value2 = 44403
$ws.Cells.Item(519,3).Value2 = 1011708
$ws.Cells.Item(519,4).Value2 = 23087
$ws.Cells.Item(519,5).Value2 = 928722
$ws.Cells.Item(519,6).Value2 = 59899
$ws.Cells.Item(520,1).Value2 = "Pakistan"
$ws.Cells.Item(520,2).Value2 = 44404
$ws.Cells.Item(520,3).Value2 = 1015827
$ws.Cells.Item(520,4).Value2 = 23133
$ws.Cells.Item(520,5).Value2 = 935742
$ws.Cells.Item(520,6).Value2 = 56952
$ws.Cells.Item(521,1).Value2 = "Pakistan"
$ws.Cells.Item(521,2).Value2 = 44405
$ws.Cells.Item(521,3).Value2 = 1020324
$ws.Cells.Item(521,4).Value2 = 23209
$ws.Cells.Item(521,5).Value2 = 937354
$ws.Cells.Item(521,6).Value2 = 59761
$ws.Cells.Item(522,1).Value2 = "Pakistan"
$ws.Cells.Item(522,2).Value2 = 44406
$ws.Cells.Item(522,3).Value2 = 1024861
$ws.Cells.Item(522,4).Value2 = 23295
$ws.Cells.Item(522,5).Value2 = 938843
$ws.Cells.Item(522,6).Value2 = 62723
$ws.Cells.Item(523,1).Value2 = "Pakistan"
$ws.Cells.Item(523,2).Value2 = 44407
$ws.Cells.Item(523,3).Value2 = 1029811
$ws.Cells.Item(523,4).Value2 = 23360
$ws.Cells.Item(523,5).Value2 = 940164
$ws.Cells.Item(523,6).Value2 = 66287
$ws.Cells.Item(524,1).Value2 = "Pakistan"
$ws.Cells.Item(524,2).Value2 = 44408
$ws.Cells.Item(524,3).Value2 = 1034837
$ws.Cells.Item(524,4).Value2 = 23422
$ws.Cells.Item(524,5).Value2 = 941659
$ws.Cells.Item(524,6).Value2 = 69756
$ws.Cells.Item(525,1).Value2 = "Pakistan"
$ws.Cells.Item(525,2).Value2 = 44409
$ws.Cells.Item(525,3).Value2 = 1039695
$ws.Cells.Item(525,4).Value2 = 23462
$ws.Cells.Item(525,5).Value2 = 943020
$ws.Cells.Item(525,6).Value2 = 73213
$ws.Cells.Item(526,1).Value2 = "Pakistan"
$ws.Cells.Item(526,2).Value2 = 44410
$ws.Cells.Item(526,3).Value2 = 1043277
$ws.Cells.Item(526,4).Value2 = 23529
$ws.Cells.Item(526,5).Value2 = 944375
$ws.Cells.Item(526,6).Value2 = 75373
$ws.Cells.Item(527,1).Value2 = "Pakistan"
$ws.Cells.Item(527,2).Value2 = 44411
$ws.Cells.Item(527,3).Value2 = 1047999
$ws.Cells.Item(527,4).Value2 = 23575
$ws.Cells.Item(527,5).Value2 = 945829
$ws.Cells.Item(527,6).Value2 = 78595
$ws.Cells.Item(528,1).Value2 = "Pakistan"
$ws.Cells.Item(528,2).Value2 = 44412
$ws.Cells.Item(528,3).Value2 = 1053660
$ws.Cells.Item(528,4).Value2 = 23635
$ws.Cells.Item(528,5).Value2 = 952616
$ws.Cells.Item(528,6).Value2 = 77409
$ws.Cells.Item(529,1).Value2 = "Pakistan"
$ws.Cells.Item(529,2).Value2 = 44413
$ws.Cells.Item(529,3).Value2 = 1058405
$ws.Cells.Item(529,4).Value2 = 23702
$ws.Cells.Item(529,5).Value2 = 954711
$ws.Cells.Item(529,6).Value2 = 1034703
$ws.Cells.Item(530,1).Value2 = "Pakistan"
$ws.Cells.Item(530,2).Value2 = 44414
$ws.Cells.Item(530,3).Value2 = 1063125
$ws.Cells.Item(530,4).Value2 = 23797
$ws.Cells.Item(530,5).Value2 = 959491
$ws.Cells.Item(530,6).Value2 = 1039328
$ws.Cells.Item(531,1).Value2 = "Pakistan"
$ws.Cells.Item(531,2).Value2 = 44415
$ws.Cells.Item(531,3).Value2 = 1067580
$ws.Cells.Item(531,4).Value2 = 23865
$ws.Cells.Item(531,5).Value2 = 961639
$ws.Cells.Item(531,6).Value2 = 1043715
$ws.Cells.Item(532,1).Value2 = "Pakistan"
$ws.Cells.Item(532,2).Value2 = 44416
$ws.Cells.Item(532,3).Value2 = 1071620
$ws.Cells.Item(532,4).Value2 = 23918
$ws.Cells.Item(532,5).Value2 = 964404
$ws.Cells.Item(532,6).Value2 = 1047702
$ws.Cells.Item(533,1).Value2 = "Pakistan"
$ws.Cells.Item(533,2).Value2 = 44417
$ws.Cells.Item(533,3).Value2 = 1075504
$ws.Cells.Item(533,4).Value2 = 24004
$ws.Cells.Item(533,5).Value2 = 967073
$ws.Cells.Item(533,6).Value2 = 1051500
$ws.Cells.Item(534,1).Value2 = "Pakistan"
$ws.Cells.Item(534,2).Value2 = 44418
$ws.Cells.Item(534,3).Value2 = 1080360
$ws.Cells.Item(534,4).Value2 = 24085
$ws.Cells.Item(534,5).Value2 = 972098
$ws.Cells.Item(534,6).Value2 = 1056275
$ws.Cells.Item(535,1).Value2 = "Pakistan"
$ws.Cells.Item(535,2).Value2 = 44419
$ws.Cells.Item(535,3).Value2 = 1085294
$ws.Cells.Item(535,4).Value2 = 24187
$ws.Cells.Item(535,5).Value2 = 975474
$ws.Cells.Item(535,6).Value2 = 1061107
$ws.Cells.Item(536,1).Value2 = "Pakistan"
$ws.Cells.Item(536,2).Value2 = 44420
$ws.Cells.Item(536,3).Value2 = 1089913
$ws.Cells.Item(536,4).Value2 = 24266
$ws.Cells.Item(536,5).Value2 = 979411
$ws.Cells.Item(536,6).Value2 = 1065647
$ws.Cells.Item(537,1).Value2 = "Pakistan"
$ws.Cells.Item(537,2).Value2 = 44421
$ws.Cells.Item(537,3).Value2 = 1094699
$ws.Cells.Item(537,4).Value2 = 24339
$ws.Cells.Item(537,5).Value2 = 983754
$ws.Cells.Item(537,6).Value2 = 1070360
$ws.Cells.Item(538,1).Value2 = "Pakistan"
$ws.Cells.Item(538,2).Value2 = 44422
$ws.Cells.Item(538,3).Value2 = 1098410
$ws.Cells.Item(538,4).Value2 = 24406
$ws.Cells.Item(538,5).Value2 = 986795
$ws.Cells.Item(538,6).Value2 = 1070360
$ws.Cells.Item(539,1).Value2 = "Pakistan"
$ws.Cells.Item(539,2).Value2 = 44423
$ws.Cells.Item(539,3).Value2 = 1102079
$ws.Cells.Item(539,4).Value2 = 24478
$ws.Cells.Item(539,5).Value2 = 989013
$ws.Cells.Item(539,6).Value2 = 1077601
$ws.Cells.Item(540,1).Value2 = "Pakistan"
$ws.Cells.Item(540,2).Value2 = 44424
$ws.Cells.Item(540,3).Value2 = 1105300
$ws.Cells.Item(540,4).Value2 = 24573
$ws.Cells.Item(540,5).Value2 = 993304
$ws.Cells.Item(540,6).Value2 = 1080727
$ws.Cells.Item(541,1).Value2 = "Pakistan"
$ws.Cells.Item(541,2).Value2 = 44425
$ws.Cells.Item(541,3).Value2 = 1109274
$ws.Cells.Item(541,4).Value2 = 24639
$ws.Cells.Item(541,5).Value2 = 996426
$ws.Cells.Item(541,6).Value2 = 1080727
$ws.Cells.Item(542,1).Value2 = "Pakistan"
$ws.Cells.Item(542,2).Value2 = 44426
$ws.Cells.Item(542,3).Value2 = 1113647
$ws.Cells.Item(542,4).Value2 = 24713
$ws.Cells.Item(542,5).Value2 = 999403
$ws.Cells.Item(542,6).Value2 = 1084635
$ws.Cells.Item(543,1).Value2 = "Pakistan"
$ws.Cells.Item(543,2).Value2 = 44427
$ws.Cells.Item(543,3).Value2 = 1116886
$ws.Cells.Item(543,4).Value2 = 24783
$ws.Cells.Item(543,5).Value2 = 1002430
$ws.Cells.Item(543,6).Value2 = 1092103
$ws.Cells.Item(544,1).Value2 = "Pakistan"
$ws.Cells.Item(544,2).Value2 = 44428
$ws.Cells.Item(544,3).Value2 = 1119970
$ws.Cells.Item(544,4).Value2 = 24848
$ws.Cells.Item(544,5).Value2 = 1006078
$ws.Cells.Item(544,6).Value2 = 1095122
$ws.Cells.Item(545,1).Value2 = "Pakistan"
$ws.Cells.Item(545,2).Value2 = 44429
$ws.Cells.Item(545,3).Value2 = 1123812
$ws.Cells.Item(545,4).Value2 = 24923
$ws.Cells.Item(545,5).Value2 = 1009555
$ws.Cells.Item(545,6).Value2 = 1098889
$ws.Cells.Item(546,1).Value2 = "Pakistan"
$ws.Cells.Item(546,2).Value2 = 44430
$ws.Cells.Item(546,3).Value2 = 1127584
$ws.Cells.Item(546,4).Value2 = 25003
$ws.Cells.Item(546,5).Value2 = 1012662
$ws.Cells.Item(546,6).Value2 = 1102581
$ws.Cells.Item(547,1).Value2 = "Pakistan"
$ws.Cells.Item(547,2).Value2 = 44431
$ws.Cells.Item(547,3).Value2 = 1131659
$ws.Cells.Item(547,4).Value2 = 25094
$ws.Cells.Item(547,5).Value2 = 1015519
$ws.Cells.Item(547,6).Value2 = 1106565

# 3) Fill the G/H/I (New Confirmed/New Deaths/New Recovered) formulas in the same
#    shared-formula groupings as the authored workbook (507:519, 520, 521, 522:534, 535:547)
$ws.Range("G507:G519").Formula = "=C507-C506"
$ws.Range("H507:H519").Formula = "=D507-D506"
$ws.Range("I507:I519").Formula = "=E507-E506"
$ws.Range("G520").Formula = "=C520-C519"
$ws.Range("H520").Formula = "=D520-D519"
$ws.Range("I520").Formula = "=E520-E519"
$ws.Range("G521").Formula = "=C521-C520"
$ws.Range("H521").Formula = "=D521-D520"
$ws.Range("I521").Formula = "=E521-E520"
$ws.Range("G522:G534").Formula = "=C522-C521"
$ws.Range("H522:H534").Formula = "=D522-D521"
$ws.Range("I522:I534").Formula = "=E522-E521"
$ws.Range("G535:G547").Formula = "=C535-C534"
$ws.Range("H535:H547").Formula = "=D535-D534"
$ws.Range("I535:I547").Formula = "=E535-E534"

# 4) Highlight the two cells that got a yellow fill in column C (large jumps flagged by the author)
$ws.Range("C522").Interior.Color = 65535
$ws.Range("C535").Interior.Color = 65535

# 5) Move the active selection to follow the new bottom of the data (the header row stays frozen)
$ws.Activate()
$ws.Range("K539").Select()

Write-Output "edit complete"
